$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in column A (Ni values for elements 4 and 5)
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# Update the active selection to match the post-edit cursor position
$ws.Range("I11").Select()
